$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 406.13043
$ws.Range("I107").Value = 310.90475
$ws.Range("K107").Value = 310.90475
$ws.Range("M107").Value = 1609.09525
$ws.Range("H127").Value = 589.1429000000001
$ws.Range("I127").Value = 480.69232
$ws.Range("J127").Value = 1999
$ws.Range("K127").Value = 1442.07696
$ws.Range("L127").Value = 5997
$ws.Range("M127").Value = 3517.92304
$ws.Range("N127").Value = -15917
$ws.Range("H129").Value = 2335.4285
$ws.Range("I129").Value = 884.2857
$ws.Range("J129").Value = 3786.5715
$ws.Range("K129").Value = 2652.8571
$ws.Range("L129").Value = 11359.7145
$ws.Range("M129").Value = 2347.1429
$ws.Range("N129").Value = -21359.7145
$ws.Range("H132").Value = 11032.275
$ws.Range("I132").Value = 6906.591
$ws.Range("K132").Value = 20719.773
$ws.Range("M132").Value = -18189.773
$ws.Range("H135").Value = 3911.6
$ws.Range("J135").Value = 5033.5
$ws.Range("L135").Value = 45301.5
$ws.Range("N135").Value = -50371.5
$ws.Range("H137").Value = 2506068.8
$ws.Range("I137").Value = 2778687.2
$ws.Range("K137").Value = 8336061.600000001
$ws.Range("M137").Value = -8333511.600000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4629.523
$ws.Range("I32").Value = 3542.4375
$ws.Range("J32").Value = 7528.4165
$ws.Range("K32").Value = 3542.4375
$ws.Range("L32").Value = 7528.4165
$ws.Range("M32").Value = -3255.4375
$ws.Range("N32").Value = -8102.4165
$ws.Range("H61").Value = 3091.8845
$ws.Range("I61").Value = 2256.5264
$ws.Range("K61").Value = 2256.5264
$ws.Range("M61").Value = -2044.5264
$ws.Range("H74").Value = 349558.2
$ws.Range("I74").Value = 506691
$ws.Range("K74").Value = 506691
$ws.Range("M74").Value = -505817
$ws.Range("H77").Value = 349558.2
$ws.Range("I77").Value = 506691
$ws.Range("K77").Value = 2533455
$ws.Range("M77").Value = -2529087
$ws.Range("H97").Value = 575
$ws.Range("I97").Value = 575
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 575
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -79
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 2997.8262
$ws.Range("I110").Value = 1665
$ws.Range("K110").Value = 1665
$ws.Range("M110").Value = 380
$ws.Range("H132").Value = 1878.279
$ws.Range("J132").Value = 4038.5
$ws.Range("L132").Value = 12115.5
$ws.Range("N132").Value = -17175.5
$ws.Range("H136").Value = 3091.8845
$ws.Range("I136").Value = 2256.5264
$ws.Range("K136").Value = 6769.5792
$ws.Range("M136").Value = -4219.5792

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 62508300
$ws.Range("I20").Value = 83342900
$ws.Range("K20").Value = 83342900
$ws.Range("M20").Value = -83342653
$ws.Range("H21").Value = 20220
$ws.Range("J21").Value = 20220
$ws.Range("L21").Value = 20220
$ws.Range("N21").Value = -20692

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3298.077
$ws.Range("I31").Value = 2067.524
$ws.Range("K31").Value = 2067.524
$ws.Range("M31").Value = -1772.524
$ws.Range("H34").Value = 3298.077
$ws.Range("I34").Value = 2067.524
$ws.Range("K34").Value = 2067.524
$ws.Range("M34").Value = -1865.524
$ws.Range("H58").Value = 2745
$ws.Range("I58").Value = 2019.9375
$ws.Range("J58").Value = 3711.75
$ws.Range("K58").Value = 2019.9375
$ws.Range("L58").Value = 3711.75
$ws.Range("M58").Value = -1816.9375
$ws.Range("N58").Value = -4117.75
$ws.Range("H134").Value = 2640.5557
$ws.Range("I134").Value = 2630.3635
$ws.Range("K134").Value = 7891.0905
$ws.Range("M134").Value = -5356.0905
$ws.Range("H136").Value = 2745
$ws.Range("I136").Value = 2019.9375
$ws.Range("J136").Value = 3711.75
$ws.Range("K136").Value = 6059.8125
$ws.Range("L136").Value = 11135.25
$ws.Range("M136").Value = -3509.8125
$ws.Range("N136").Value = -16235.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2470065
$ws.Range("I4").Value = 2150308.5
$ws.Range("K4").Value = 6450925.5
$ws.Range("M4").Value = -6450813.5
$ws.Range("H7").Value = 550
$ws.Range("I7").Value = 550
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1650
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1538
$ws.Range("N7").ClearContents()
$ws.Range("H11").Value = 419.8
$ws.Range("I11").Value = 266.66666
$ws.Range("J11").Value = 649.5
$ws.Range("K11").Value = 799.9999799999999
$ws.Range("L11").Value = 1948.5
$ws.Range("M11").Value = -659.9999799999999
$ws.Range("N11").Value = -2228.5
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("K13").Value = 300
$ws.Range("M13").Value = -132
$ws.Range("H21").Value = 6538721.5
$ws.Range("I21").Value = 13891014
$ws.Range("J21").Value = 3350.2222
$ws.Range("K21").Value = 41673042
$ws.Range("L21").Value = 10050.6666
$ws.Range("M21").Value = -41672869
$ws.Range("N21").Value = -10396.6666
$ws.Range("H24").Value = 3068.8
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 3068.8
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 9206.400000000001
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -9666.400000000001
$ws.Range("H25").Value = 6590.4546
$ws.Range("I25").Value = 999.5
$ws.Range("J25").Value = 7832.8887
$ws.Range("K25").Value = 2998.5
$ws.Range("L25").Value = 23498.6661
$ws.Range("M25").Value = -2829.5
$ws.Range("N25").Value = -23836.6661
$ws.Range("H29").Value = 754.875
$ws.Range("I29").Value = 96.666664
$ws.Range("J29").Value = 1149.8
$ws.Range("K29").Value = 289.999992
$ws.Range("L29").Value = 3449.4
$ws.Range("M29").Value = -12.99999200000002
$ws.Range("N29").Value = -4003.4
$ws.Range("H30").Value = 6590.4546
$ws.Range("I30").Value = 999.5
$ws.Range("J30").Value = 7832.8887
$ws.Range("K30").Value = 2998.5
$ws.Range("L30").Value = 23498.6661
$ws.Range("M30").Value = -2896.5
$ws.Range("N30").Value = -23702.6661
$ws.Range("H36").Value = 2374.5
$ws.Range("I36").Value = 899
$ws.Range("J36").Value = 2866.3333
$ws.Range("K36").Value = 2697
$ws.Range("L36").Value = 8598.999899999999
$ws.Range("M36").Value = -2528
$ws.Range("N36").Value = -8936.999899999999
$ws.Range("H39").Value = 6079.5557
$ws.Range("J39").Value = 6326.1177
$ws.Range("L39").Value = 18978.3531
$ws.Range("N39").Value = -19566.3531
$ws.Range("H43").Value = 3033.3333
$ws.Range("J43").Value = 3033.3333
$ws.Range("L43").Value = 9099.999899999999
$ws.Range("N43").Value = -9327.999899999999
$ws.Range("H51").Value = 350
$ws.Range("I51").Value = 350
$ws.Range("K51").Value = 1050
$ws.Range("M51").Value = -590
$ws.Range("H58").Value = 4248.75
$ws.Range("I58").Value = 3998
$ws.Range("J58").Value = 4332.3335
$ws.Range("K58").Value = 11994
$ws.Range("L58").Value = 12997.0005
$ws.Range("M58").Value = -11866
$ws.Range("N58").Value = -13253.0005
$ws.Range("H62").Value = 8299.833000000001
$ws.Range("J62").Value = 8299.833000000001
$ws.Range("L62").Value = 24899.499
$ws.Range("N62").Value = -26271.499
$ws.Range("H65").Value = 8299.833000000001
$ws.Range("J65").Value = 8299.833000000001
$ws.Range("L65").Value = 74698.497
$ws.Range("N65").Value = -81562.497

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27266.25
$ws.Range("J46").Value = 36021.668
$ws.Range("L46").Value = 36021.668
$ws.Range("N46").Value = -36333.668
$ws.Range("H57").Value = 64933.715
$ws.Range("J57").Value = 64933.715
$ws.Range("L57").Value = 64933.715
$ws.Range("N57").Value = -66573.715
$ws.Range("H70").Value = 226500.33
$ws.Range("I70").Value = 254062.88
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 254062.88
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -253792.88
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 226500.33
$ws.Range("I73").Value = 254062.88
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 254062.88
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -253126.88
$ws.Range("N73").Value = -7872
$ws.Range("H132").Value = 3783.2942
$ws.Range("I132").Value = 3301.7273
$ws.Range("J132").Value = 4666.1665
$ws.Range("K132").Value = 9905.1819
$ws.Range("L132").Value = 13998.4995
$ws.Range("M132").Value = -7375.1819
$ws.Range("N132").Value = -19058.4995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9150.200000000001
$ws.Range("I40").Value = 8000
$ws.Range("J40").Value = 9437.75
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 9437.75
$ws.Range("M40").Value = -7864
$ws.Range("N40").Value = -9709.75
$ws.Range("H122").Value = 3390.3333
$ws.Range("J122").Value = 3176.875
$ws.Range("L122").Value = 9530.625
$ws.Range("N122").Value = -14430.625
$ws.Range("H132").Value = 9710.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 9710.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 29131.2
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -34191.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13890343
$ws.Range("I122").Value = 1344.2307
$ws.Range("K122").Value = 4032.6921
$ws.Range("M122").Value = -1582.6921
$ws.Range("H132").Value = 3425.7715
$ws.Range("I132").Value = 3662.9644
$ws.Range("K132").Value = 10988.8932
$ws.Range("M132").Value = -8458.893199999999
$ws.Range("H136").Value = 6012
$ws.Range("I136").Value = 2013.9
$ws.Range("J136").Value = 2013.9
$ws.Range("K136").Value = 6041.700000000001
$ws.Range("M136").Value = -3491.700000000001
